$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the "Brief Description of Project and Methods" answer.
#    The empty paragraph right after the "Brief Description..." label gets a
#    new sibling paragraph (bold-complex-script run formatting, 3 runs) with
#    the project description text, and the "_GoBack" bookmark (which used to
#    sit at the end of the "Date:" answer further down) now lives here -
#    this mirrors Word's behaviour of moving _GoBack to the most recently
#    edited spot.
# ---------------------------------------------------------------------------

$labelIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Brief Description of Project and Methods:") {
        $labelIndex = $i
        break
    }
}

if ($labelIndex -eq -1) {
    throw "Could not locate the 'Brief Description of Project and Methods:' paragraph"
}

$emptyPara = $d.Paragraphs.Item($labelIndex + 1)
$insertPoint = $d.Range($emptyPara.Range.End, $emptyPara.Range.End)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:bCs/></w:rPr><w:t>I will conduct an analysis how accessibility to cities is affecting biodiversity globally. I will mainly do data manipulation and analysis using the Google Earth Engine and R</w:t></w:r>' +
  '<w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> on my personal computer</w:t></w:r>' +
  '<w:r><w:rPr><w:bCs/></w:rPr><w:t>. There is no field work involved in my research.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'

$insertPoint.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 2) Remove the stale "_GoBack" bookmark from the "Date: 14/01/20" answer
#    (the student-signature date cell). It has now been superseded by the
#    bookmark added above, so only one "_GoBack" should remain in the doc.
# ---------------------------------------------------------------------------

$dateIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "14/01/20") {
        $dateIndex = $i
        break
    }
}

if ($dateIndex -eq -1) {
    throw "Could not locate the 'Date: 14/01/20' paragraph"
}

$datePara = $d.Paragraphs.Item($dateIndex)
$fullRange = $datePara.Range

$dateParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
  'w14:paraId="7CAED0F0" w14:textId="779D4821" w:rsidR="00824EEF" w:rsidRDefault="0084359B">' +
  '<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>14/01/20</w:t></w:r>' +
  '</w:p>'

$fullRange.InsertXML($dateParaXml)

Write-Output "edit complete"
